# Update the date heading and every arithmetic answer in the practice-sheet
# table. Each cell's text is unique within the document, so a plain
# Find/Replace (MatchWholeWord, no wildcards) targeting each old value in
# document order is sufficient and unambiguous.
$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-10-23 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-24 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("42-23=19", $true, $false, $false, $false, $false, $true, 1, $false, "70-4=66", 2) | Out-Null
$d.Content.Find.Execute("36+30=66", $true, $false, $false, $false, $false, $true, 1, $false, "89-57=32", 2) | Out-Null
$d.Content.Find.Execute("91-20=71", $true, $false, $false, $false, $false, $true, 1, $false, "92+7=99", 2) | Out-Null
$d.Content.Find.Execute("68-58=10", $true, $false, $false, $false, $false, $true, 1, $false, "3+74=77", 2) | Out-Null
$d.Content.Find.Execute("70-64=6", $true, $false, $false, $false, $false, $true, 1, $false, "15+72=87", 2) | Out-Null
$d.Content.Find.Execute("45+11=56", $true, $false, $false, $false, $false, $true, 1, $false, "65+15=80", 2) | Out-Null
$d.Content.Find.Execute("76+9=85", $true, $false, $false, $false, $false, $true, 1, $false, "57-25=32", 2) | Out-Null
$d.Content.Find.Execute("92-47=45", $true, $false, $false, $false, $false, $true, 1, $false, "15+21=36", 2) | Out-Null
$d.Content.Find.Execute("21+36=57", $true, $false, $false, $false, $false, $true, 1, $false, "15+28=43", 2) | Out-Null
$d.Content.Find.Execute("48+50=98", $true, $false, $false, $false, $false, $true, 1, $false, "57+24=81", 2) | Out-Null
$d.Content.Find.Execute("98-22=76", $true, $false, $false, $false, $false, $true, 1, $false, "0+74=74", 2) | Out-Null
$d.Content.Find.Execute("89-46=43", $true, $false, $false, $false, $false, $true, 1, $false, "75-64=11", 2) | Out-Null
$d.Content.Find.Execute("50-25=25", $true, $false, $false, $false, $false, $true, 1, $false, "95-82=13", 2) | Out-Null
$d.Content.Find.Execute("34+64=98", $true, $false, $false, $false, $false, $true, 1, $false, "59-23=36", 2) | Out-Null
$d.Content.Find.Execute("8+34=42", $true, $false, $false, $false, $false, $true, 1, $false, "43+4=47", 2) | Out-Null
$d.Content.Find.Execute("63+6=69", $true, $false, $false, $false, $false, $true, 1, $false, "99-53=46", 2) | Out-Null
$d.Content.Find.Execute("25+68=93", $true, $false, $false, $false, $false, $true, 1, $false, "46-4=42", 2) | Out-Null
$d.Content.Find.Execute("52-15=37", $true, $false, $false, $false, $false, $true, 1, $false, "34+18=52", 2) | Out-Null
$d.Content.Find.Execute("14+64=78", $true, $false, $false, $false, $false, $true, 1, $false, "26+44=70", 2) | Out-Null
$d.Content.Find.Execute("7+65=72", $true, $false, $false, $false, $false, $true, 1, $false, "44+13=57", 2) | Out-Null
$d.Content.Find.Execute("12+17=29", $true, $false, $false, $false, $false, $true, 1, $false, "67-31=36", 2) | Out-Null
$d.Content.Find.Execute("52+19=71", $true, $false, $false, $false, $false, $true, 1, $false, "16+13=29", 2) | Out-Null
$d.Content.Find.Execute("36-5=31", $true, $false, $false, $false, $false, $true, 1, $false, "96-9=87", 2) | Out-Null
$d.Content.Find.Execute("14+80=94", $true, $false, $false, $false, $false, $true, 1, $false, "94-57=37", 2) | Out-Null
$d.Content.Find.Execute("96-68=28", $true, $false, $false, $false, $false, $true, 1, $false, "44+5=49", 2) | Out-Null
$d.Content.Find.Execute("20-12=8", $true, $false, $false, $false, $false, $true, 1, $false, "94-57=37", 2) | Out-Null
$d.Content.Find.Execute("81+11=92", $true, $false, $false, $false, $false, $true, 1, $false, "56-41=15", 2) | Out-Null
$d.Content.Find.Execute("87-78=9", $true, $false, $false, $false, $false, $true, 1, $false, "85-77=8", 2) | Out-Null
$d.Content.Find.Execute("33-26=7", $true, $false, $false, $false, $false, $true, 1, $false, "5+35=40", 2) | Out-Null
$d.Content.Find.Execute("59+0=59", $true, $false, $false, $false, $false, $true, 1, $false, "41+14=55", 2) | Out-Null
$d.Content.Find.Execute("89+2=91", $true, $false, $false, $false, $false, $true, 1, $false, "85-63=22", 2) | Out-Null
$d.Content.Find.Execute("50+12=62", $true, $false, $false, $false, $false, $true, 1, $false, "56-42=14", 2) | Out-Null
$d.Content.Find.Execute("38+29=67", $true, $false, $false, $false, $false, $true, 1, $false, "71-64=7", 2) | Out-Null
$d.Content.Find.Execute("90-30=60", $true, $false, $false, $false, $false, $true, 1, $false, "40+44=84", 2) | Out-Null
$d.Content.Find.Execute("92-91=1", $true, $false, $false, $false, $false, $true, 1, $false, "24+38=62", 2) | Out-Null
$d.Content.Find.Execute("5+15=20", $true, $false, $false, $false, $false, $true, 1, $false, "15+22=37", 2) | Out-Null
$d.Content.Find.Execute("3+11=14", $true, $false, $false, $false, $false, $true, 1, $false, "33+19=52", 2) | Out-Null
$d.Content.Find.Execute("13+19=32", $true, $false, $false, $false, $false, $true, 1, $false, "23+21=44", 2) | Out-Null
$d.Content.Find.Execute("58+0=58", $true, $false, $false, $false, $false, $true, 1, $false, "6+46=52", 2) | Out-Null
$d.Content.Find.Execute("85-55=30", $true, $false, $false, $false, $false, $true, 1, $false, "74-60=14", 2) | Out-Null
$d.Content.Find.Execute("44+32=76", $true, $false, $false, $false, $false, $true, 1, $false, "48+21=69", 2) | Out-Null
$d.Content.Find.Execute("49+28=77", $true, $false, $false, $false, $false, $true, 1, $false, "47-35=12", 2) | Out-Null
$d.Content.Find.Execute("62+10=72", $true, $false, $false, $false, $false, $true, 1, $false, "30-19=11", 2) | Out-Null
$d.Content.Find.Execute("22+44=66", $true, $false, $false, $false, $false, $true, 1, $false, "83+5=88", 2) | Out-Null
$d.Content.Find.Execute("29+23=52", $true, $false, $false, $false, $false, $true, 1, $false, "21+71=92", 2) | Out-Null
$d.Content.Find.Execute("88-50=38", $true, $false, $false, $false, $false, $true, 1, $false, "42-24=18", 2) | Out-Null
$d.Content.Find.Execute("70-27=43", $true, $false, $false, $false, $false, $true, 1, $false, "43+19=62", 2) | Out-Null
$d.Content.Find.Execute("24+14=38", $true, $false, $false, $false, $false, $true, 1, $false, "16+41=57", 2) | Out-Null
$d.Content.Find.Execute("43-17=26", $true, $false, $false, $false, $false, $true, 1, $false, "79-39=40", 2) | Out-Null
$d.Content.Find.Execute("78-75=3", $true, $false, $false, $false, $false, $true, 1, $false, "94-37=57", 2) | Out-Null
$d.Content.Find.Execute("11+59=70", $true, $false, $false, $false, $false, $true, 1, $false, "59-17=42", 2) | Out-Null
$d.Content.Find.Execute("5+85=90", $true, $false, $false, $false, $false, $true, 1, $false, "29+8=37", 2) | Out-Null
$d.Content.Find.Execute("17+3=20", $true, $false, $false, $false, $false, $true, 1, $false, "71+23=94", 2) | Out-Null
$d.Content.Find.Execute("91-62=29", $true, $false, $false, $false, $false, $true, 1, $false, "17+70=87", 2) | Out-Null
$d.Content.Find.Execute("83+14=97", $true, $false, $false, $false, $false, $true, 1, $false, "71+23=94", 2) | Out-Null
$d.Content.Find.Execute("93-16=77", $true, $false, $false, $false, $false, $true, 1, $false, "50-15=35", 2) | Out-Null
$d.Content.Find.Execute("48-39=9", $true, $false, $false, $false, $false, $true, 1, $false, "89-36=53", 2) | Out-Null
$d.Content.Find.Execute("89-81=8", $true, $false, $false, $false, $false, $true, 1, $false, "80-31=49", 2) | Out-Null
$d.Content.Find.Execute("66+27=93", $true, $false, $false, $false, $false, $true, 1, $false, "95-93=2", 2) | Out-Null
$d.Content.Find.Execute("31+45=76", $true, $false, $false, $false, $false, $true, 1, $false, "92-55=37", 2) | Out-Null
$d.Content.Find.Execute("70-32=38", $true, $false, $false, $false, $false, $true, 1, $false, "8+50=58", 2) | Out-Null
$d.Content.Find.Execute("9+35=44", $true, $false, $false, $false, $false, $true, 1, $false, "60-51=9", 2) | Out-Null
$d.Content.Find.Execute("66+1=67", $true, $false, $false, $false, $false, $true, 1, $false, "32-10=22", 2) | Out-Null
$d.Content.Find.Execute("82-47=35", $true, $false, $false, $false, $false, $true, 1, $false, "5+19=24", 2) | Out-Null
$d.Content.Find.Execute("85-61=24", $true, $false, $false, $false, $false, $true, 1, $false, "76+23=99", 2) | Out-Null
$d.Content.Find.Execute("89-54=35", $true, $false, $false, $false, $false, $true, 1, $false, "76-47=29", 2) | Out-Null
$d.Content.Find.Execute("75-75=0", $true, $false, $false, $false, $false, $true, 1, $false, "6+90=96", 2) | Out-Null
$d.Content.Find.Execute("91-79=12", $true, $false, $false, $false, $false, $true, 1, $false, "4+53=57", 2) | Out-Null
$d.Content.Find.Execute("1+17=18", $true, $false, $false, $false, $false, $true, 1, $false, "68+31=99", 2) | Out-Null
$d.Content.Find.Execute("20+43=63", $true, $false, $false, $false, $false, $true, 1, $false, "89-87=2", 2) | Out-Null
$d.Content.Find.Execute("17+10=27", $true, $false, $false, $false, $false, $true, 1, $false, "76-67=9", 2) | Out-Null
$d.Content.Find.Execute("60+16=76", $true, $false, $false, $false, $false, $true, 1, $false, "32+26=58", 2) | Out-Null
$d.Content.Find.Execute("49-47=2", $true, $false, $false, $false, $false, $true, 1, $false, "48-6=42", 2) | Out-Null
$d.Content.Find.Execute("86-39=47", $true, $false, $false, $false, $false, $true, 1, $false, "68-23=45", 2) | Out-Null
$d.Content.Find.Execute("83-65=18", $true, $false, $false, $false, $false, $true, 1, $false, "35+26=61", 2) | Out-Null
$d.Content.Find.Execute("7+46=53", $true, $false, $false, $false, $false, $true, 1, $false, "1+54=55", 2) | Out-Null
$d.Content.Find.Execute("66+10=76", $true, $false, $false, $false, $false, $true, 1, $false, "94-7=87", 2) | Out-Null
$d.Content.Find.Execute("64+1=65", $true, $false, $false, $false, $false, $true, 1, $false, "18+14=32", 2) | Out-Null
$d.Content.Find.Execute("44-3=41", $true, $false, $false, $false, $false, $true, 1, $false, "53+33=86", 2) | Out-Null
$d.Content.Find.Execute("71-29=42", $true, $false, $false, $false, $false, $true, 1, $false, "55-43=12", 2) | Out-Null
$d.Content.Find.Execute("17+29=46", $true, $false, $false, $false, $false, $true, 1, $false, "33-8=25", 2) | Out-Null
$d.Content.Find.Execute("28+53=81", $true, $false, $false, $false, $false, $true, 1, $false, "66+16=82", 2) | Out-Null
$d.Content.Find.Execute("93-69=24", $true, $false, $false, $false, $false, $true, 1, $false, "80-14=66", 2) | Out-Null
$d.Content.Find.Execute("85-59=26", $true, $false, $false, $false, $false, $true, 1, $false, "6+23=29", 2) | Out-Null
$d.Content.Find.Execute("73-68=5", $true, $false, $false, $false, $false, $true, 1, $false, "40-39=1", 2) | Out-Null
$d.Content.Find.Execute("57+12=69", $true, $false, $false, $false, $false, $true, 1, $false, "76-47=29", 2) | Out-Null
$d.Content.Find.Execute("25-12=13", $true, $false, $false, $false, $false, $true, 1, $false, "39-31=8", 2) | Out-Null
$d.Content.Find.Execute("83-82=1", $true, $false, $false, $false, $false, $true, 1, $false, "70+16=86", 2) | Out-Null
$d.Content.Find.Execute("71-48=23", $true, $false, $false, $false, $false, $true, 1, $false, "7+54=61", 2) | Out-Null
$d.Content.Find.Execute("83-40=43", $true, $false, $false, $false, $false, $true, 1, $false, "11-5=6", 2) | Out-Null
$d.Content.Find.Execute("85-71=14", $true, $false, $false, $false, $false, $true, 1, $false, "95-29=66", 2) | Out-Null
$d.Content.Find.Execute("58-11=47", $true, $false, $false, $false, $false, $true, 1, $false, "70-38=32", 2) | Out-Null
$d.Content.Find.Execute("0+34=34", $true, $false, $false, $false, $false, $true, 1, $false, "58+13=71", 2) | Out-Null
$d.Content.Find.Execute("85-50=35", $true, $false, $false, $false, $false, $true, 1, $false, "75-29=46", 2) | Out-Null
$d.Content.Find.Execute("93-31=62", $true, $false, $false, $false, $false, $true, 1, $false, "72-69=3", 2) | Out-Null
$d.Content.Find.Execute("14+9=23", $true, $false, $false, $false, $false, $true, 1, $false, "42+44=86", 2) | Out-Null
$d.Content.Find.Execute("19+17=36", $true, $false, $false, $false, $false, $true, 1, $false, "48-11=37", 2) | Out-Null
$d.Content.Find.Execute("32+63=95", $true, $false, $false, $false, $false, $true, 1, $false, "2+19=21", 2) | Out-Null
$d.Content.Find.Execute("14+49=63", $true, $false, $false, $false, $false, $true, 1, $false, "85-34=51", 2) | Out-Null
$d.Content.Find.Execute("39+54=93", $true, $false, $false, $false, $false, $true, 1, $false, "32+39=71", 2) | Out-Null
